$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trigrams (column B) and frequencies (column C) for rows 2-81,
# reflecting removal of garbled contraction entries and the resulting
# upward shift, with new trigram rows appended at the bottom.
$trigrams = @(
    @('(''sea'', ''level'', ''rise'')', 13),
    @('(''reddit'', ''kotakuinact'', ''comment'')', 13),
    @('(''climat'', ''chang'', ''polit'')', 6),
    @('(''climat'', ''chang'', ''real'')', 6),
    @('(''giss'', ''nasa'', ''gov'')', 5),
    @('(''believ'', ''climat'', ''chang'')', 5),
    @('(''unit'', ''nation'', ''publish'')', 4),
    @('(''nation'', ''publish'', ''report'')', 4),
    @('(''publish'', ''report'', ''last'')', 4),
    @('(''report'', ''last'', ''month'')', 4),
    @('(''last'', ''month'', ''warn'')', 4),
    @('(''month'', ''warn'', ''drastic'')', 4),
    @('(''warn'', ''drastic'', ''refor'')', 4),
    @('(''drastic'', ''refor'', ''effort'')', 4),
    @('(''refor'', ''effort'', ''reduct'')', 4),
    @('(''effort'', ''reduct'', ''carbon'')', 4),
    @('(''reduct'', ''carbon'', ''pollut'')', 4),
    @('(''carbon'', ''pollut'', ''planet'')', 4),
    @('(''pollut'', ''planet'', ''soon'')', 4),
    @('(''planet'', ''soon'', ''face'')', 4),
    @('(''soon'', ''face'', ''irrever'')', 4),
    @('(''face'', ''irrever'', ''global'')', 4),
    @('(''irrever'', ''global'', ''catastroph'')', 4),
    @('(''web'', ''archiv'', ''web'')', 4),
    @('(''due'', ''climat'', ''chang'')', 4),
    @('(''climat'', ''chang'', ''happen'')', 4),
    @('(''climat'', ''chang'', ''go'')', 4),
    @('(''nasa'', ''gov'', ''tmp'')', 4),
    @('(''gov'', ''tmp'', ''gistemp'')', 4),
    @('(''tmp'', ''gistemp'', ''custom'')', 4),
    @('(''gistemp'', ''custom'', ''graph'')', 4),
    @('(''ice'', ''cap'', ''melt'')', 4),
    @('(''peopl'', ''3rd'', ''world'')', 4),
    @('(''unit'', ''state'', ''recov'')', 3),
    @('(''histor'', ''wildfir'', ''california'')', 3),
    @('(''climat'', ''chang'', ''mechan'')', 3),
    @('(''global'', ''warm'', ''game'')', 3),
    @('(''esrl'', ''noaa'', ''gov'')', 3),
    @('(''mbrol67'', ''reddit'', ''kotakuinact'')', 3),
    @('(''scienc'', ''sciencemag'', ''content'')', 3),
    @('(''kotakuinact'', ''comment'', ''ethicssocju'')', 3),
    @('(''comment'', ''ethicssocju'', ''eurogam'')', 3),
    @('(''ethicssocju'', ''eurogam'', ''civil'')', 3),
    @('(''eurogam'', ''civil'', ''gather'')', 3),
    @('(''civil'', ''gather'', ''ea7hvkw'')', 3),
    @('(''climat'', ''chang'', ''alway'')', 3),
    @('(''chang'', ''global'', ''warm'')', 3),
    @('(''climat'', ''chang'', ''make'')', 3),
    @('(''effect'', ''climat'', ''chang'')', 3),
    @('(''say'', ''climat'', ''chang'')', 3),
    @('(''chang'', ''polit'', ''statement'')', 3),
    @('(''climat'', ''chang'', ''actual'')', 3),
    @('(''nuclear'', ''power'', ''plant'')', 3),
    @('(''stop'', ''global'', ''warm'')', 3),
    @('(''antarct'', ''ice'', ''sheet'')', 3),
    @('(''believ'', ''man'', ''make'')', 3),
    @('(''graph'', ''png'', ''data'')', 3),
    @('(''png'', ''data'', ''giss'')', 3),
    @('(''data'', ''giss'', ''nasa'')', 3),
    @('(''climat'', ''nasa'', ''gov'')', 3),
    @('(''climat'', ''chang'', ''game'')', 3),
    @('(''climat'', ''chang'', ''understand'')', 3),
    @('(''climat'', ''chang'', ''scienc'')', 3),
    @('(''nation'', ''secur'', ''threat'')', 3),
    @('(''vike'', ''grow'', ''wheat'')', 3),
    @('(''absorb'', ''low'', ''frequenc'')', 3),
    @('(''3rd'', ''world'', ''west'')', 3),
    @('(''climat'', ''chang'', ''propos'')', 3),
    @('(''import'', ''peopl'', ''3rd'')', 3),
    @('(''tell'', ''3rd'', ''world'')', 3),
    @('(''3rd'', ''world'', ''breed'')', 3),
    @('(''world'', ''breed'', ''rabbit'')', 3),
    @('(''low'', ''iq'', ''peopl'')', 3),
    @('(''economi'', ''climat'', ''chang'')', 3),
    @('(''dislik'', ''big'', ''govern'')', 3),
    @('(''use'', ''coal'', ''oil'')', 3),
    @('(''special'', ''interest'', ''donor'')', 3),
    @('(''peopl'', ''think'', ''idiot'')', 3),
    @('(''think'', ''idiot'', ''take'')', 3),
    @('(''idiot'', ''take'', ''much'')', 3)
)

for ($i = 0; $i -lt $trigrams.Length; $i++) {
    $row = 2 + $i
    $ws.Range("B$row").Value = $trigrams[$i][0]
    $ws.Range("C$row").Value = $trigrams[$i][1]
}

